$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (the "Förändrad" date column) for rows 2-9
# from serial date 46070 (2026-02-17) to 46072 (2026-02-19)
$ws.Range("C2:C9").Value = 46072
